$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data (filtering save games) for ortiz_luis.
# Columns: B=TB, C=d2S, D=K, E=IP, F=Win (unchanged), G=sum (=B+C+D+E)

$data = @{
    2 = @{ B = 1.505614041169197;    C = 1.65323645889881;    D = 0.7127328510149897; E = 0.4998867070740569;  G = 4.371470058157054 }
    3 = @{ B = 0.06328177979961902;  C = 0.3375848360084654;  D = 0.7127328510149897; E = 6.48142807727062;     G = 7.595027544093695 }
    4 = @{ B = 0.00006486019690155054; C = 10990084.13351303;  D = 0.1529057820181812; E = 246.9852506941017;    G = 10990331.27173437 }
    5 = @{ B = 0.7287194209349384;   C = 1.65323645889881;    D = 0.7127328510149897; E = 0.4998867070740569;  G = 3.594575437922795 }
    6 = @{ B = 0.1554434735375247;   C = 0.3375848360084654;  D = 0.7127328510149897; E = 6.48142807727062;     G = 7.6871892378316 }
    7 = @{ B = 1.505614041169197;    C = 1.65323645889881;    D = 3.082599426703578;  E = 0.4998867070740569;  G = 6.741336633845642 }
    8 = @{ B = 3.182878228561681;    C = 1.65323645889881;    D = 3.082599426703578;  E = 0.4998867070740569;  G = 8.418600821238126 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
